$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (C5:C41): re-apply Text number format (style dedups to the
#     quotePrefix+Text xf that already exists in the workbook) ---
$ws.Range("C5:C41").NumberFormat = "@"

# --- Column G (G5:G40): results were numeric 30, change to the text-typed
#     results of the concatenation expressions. "10.020" and "1020.0" are
#     brand-new shared strings, so seed them first (in this order) to get
#     the same shared-string ids the original edit produced, then fill in
#     the rest of the column in sheet order. Quote-prefixed so the
#     numeric-looking text stays text. ---
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "'10.020"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "'1020.0"

$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "'1020"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "'1020"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "'1020"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "'1020"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "'1020.0"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "'1020"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "'1020"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "'1020"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "'1020"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "'10.020"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "'1020"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "'1020"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "'1020"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "'1020"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "'1020.0"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "'1020.0"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "'1020"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "'1020"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "'1020"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "'1020"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "'10.020"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "'10.020"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "'1020"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "'1020"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "'1020"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "'1020"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "'1020.0"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "'1020.0"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "'1020"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "'1020"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "'1020"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "'1020"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "'10.020"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "'10.020"

# --- Move the active selection to L14 (matches the author's final cursor
#     position when they saved) ---
$ws.Range("L14").Select()
